# Automatische test-sync: 2025-06-19 21:46:50
#
# Adds a new incoming-mail log entry to the "Logs" sheet (row 29) and
# refreshes the category tally on the "Dashboard" sheet to reflect the
# change: "IT / Technisch probleem" now has 3 occurrences (was 2) and
# swaps ranking order with "Afmelding / Nieuwsbrief" (still 3).

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Append the new mail log row (row 29) to the Logs sheet ---
$row = 29
$logs.Cells.Item($row, 1).Value = "Probleem met inloggen"
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value = "Ik kan niet inloggen op mijn account. Kunnen jullie dit oplossen?"
$logs.Cells.Item($row, 4).Value = "IT / Technisch probleem"
# Column E (Antwoord) intentionally left blank - mail not yet answered
$logs.Cells.Item($row, 6).Value = "2025-06-19 21:46:31"
$logs.Cells.Item($row, 7).Value = "Nee"

# --- Update the Dashboard category counts / ordering ---
# Row 5 becomes "IT / Technisch probleem" with count 3
$dashboard.Cells.Item(5, 1).Value = "IT / Technisch probleem"
$dashboard.Cells.Item(5, 2).Value = 3

# Row 6 becomes "Afmelding / Nieuwsbrief" with count 3
$dashboard.Cells.Item(6, 1).Value = "Afmelding / Nieuwsbrief"
$dashboard.Cells.Item(6, 2).Value = 3

# --- Grow the conditional formatting ranges on the Logs sheet to include row 29 ---
$dConditions = $logs.Range("D2:D28").FormatConditions
for ($i = 1; $i -le $dConditions.Count; $i++) {
    $dConditions.Item($i).ModifyAppliesToRange($logs.Range("D2:D29"))
}

$gConditions = $logs.Range("G2:G28").FormatConditions
for ($i = 1; $i -le $gConditions.Count; $i++) {
    $gConditions.Item($i).ModifyAppliesToRange($logs.Range("G2:G29"))
}
